$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the now-unused duplicate chart defined names (_xlchart.v1.2 / v1.3)
$wb.Names.Item("_xlchart.v1.2").Delete()
$wb.Names.Item("_xlchart.v1.3").Delete()

# New "Std" / "Relative std" header labels (bold, matching the other stat headers)
$ws.Range("D15").Value = "Std"
$ws.Range("D15").Font.Bold = $true
$ws.Range("E15").Value = "Relative std"
$ws.Range("E15").Font.Bold = $true

# New standard-deviation + relative-standard-deviation measures
$ws.Range("D16").Formula = "=STDEV(B2:B31)"
$ws.Range("E16").Formula = "= (D16 / E4) * 100"

# Move the active cell selection from M20 to P20
$ws.Range("P20").Select()
